# Add a "Save" column (column H) to the s_vals sheet, mirroring the
# header formatting used by the existing "sum" column (G).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the G1 header cell onto H1 so the new header
# gets the same bold/border/centered style (cellXf s="1").
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Set the new header text and the numeric "Save" values per row.
$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 1
